# App layout and styling
#
# Adds two new "Minus Sign" shapes (mathMinus preset geometry) to slide 2,
# matching the existing "Minus Sign 27" shape's style/appearance but at new
# positions/sizes. These are created by duplicating the existing shape so
# that the line/fill/effect/font style refs and text body come along for
# free, then repositioning, resizing and renaming the duplicates.

$pres = $ppt.ActivePresentation
$slide = $pres.Slides.Item(2)

# Best-effort: touch the presentation-level Guides collection. (No-op in
# some hosts; harmless either way.)
try {
    [void]$pres.Guides
} catch {
}

# Locate the existing "Minus Sign 27" shape to use as a style template.
$template = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.Name -eq "Minus Sign 27") {
        $template = $candidate
    }
}

if ($template -eq $null) {
    # Fallback: build a fresh mathMinus shape from scratch with the same
    # appearance if the expected template shape could not be found.
    $template = $slide.Shapes.AddShape("mathMinus", 0, 0, 23.5496454292126, 26.023345956614175)
}

# First new shape: "Minus Sign 1"
$dup1 = $template.Duplicate()
$shape1 = $dup1.Item(1)
$shape1.Name = "Minus Sign 1"
$shape1.Left = 671.8692626984251
$shape1.Top = 371.5609893818898
$shape1.Width = 23.5496454292126
$shape1.Height = 26.023345956614175
$shape1.Adjustments.Item(1) = 0.15972

# Second new shape: "Minus Sign 2"
$dup2 = $template.Duplicate()
$shape2 = $dup2.Item(1)
$shape2.Name = "Minus Sign 2"
$shape2.Left = 356.3139038077165
$shape2.Top = 357.6975250249606
$shape2.Width = 23.5496454292126
$shape2.Height = 26.023345956614175
$shape2.Adjustments.Item(1) = 0.15972

Write-Output "Added shapes: $($shape1.Name) @ ($($shape1.Left),$($shape1.Top)) and $($shape2.Name) @ ($($shape2.Left),$($shape2.Top))"
